# Revised report sections on data scaling and model architecture
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Remove the "Quitar" helper column (G) ---
# G1 held the shared string "Quitar"; G2:G9 held "=C#-C(row+10)" formulas.
$ws.Range("G1:G9").ClearContents()

# --- Remove the old MIN() summary row (row 12) ---
$ws.Range("C12").ClearContents()

# --- Add the new totals row (row 11) ---
$ws.Range("C11").Formula = "=SUM(C2:C9)"
$ws.Range("E11").Formula = "=81*8"
$ws.Range("D11").Formula = "=C11-E11"

# Apply the same number-format style (style index 1, numFmtId 1) used by
# columns D/E elsewhere in the sheet to the new D11/E11 cells.
$ws.Range("D11").NumberFormat = "0"
$ws.Range("E11").NumberFormat = "0"

# --- Drop the conditional formatting that referenced the old C2:C12 range ---
# (Re-apply an equivalent duplicate-values rule over the new data range first,
# matching the sequence of edits the author actually performed, then clear
# all conditional formatting from the sheet.)
$fc = $ws.Range("C2:C11").FormatConditions.AddUniqueValues()
$fc.DupeUnique = 1
$fc.Font.Color = -16383844
$fc.Interior.Color = 13551615
$ws.Cells.FormatConditions.Delete()

# --- Update the active selection to reflect the new working cell ---
$ws.Range("D11").Select()

$wb.Save()
